$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.171352624893188
$ws.Range("B1").Value = 4.763513565063477
$ws.Range("C1").Value = 3.228675603866577
$ws.Range("D1").Value = 1.700800776481628
$ws.Range("E1").Value = 1.507060885429382
